$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Trim the sentence that references the old "standard form
#    certificate" wording so the paragraph ends right after the domain
#    placeholder.
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*{{{domain}}}*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $d.Range($target.Range.Start, $target.Range.End)
    $openQuote = [char]0x201C
    $closeQuote = [char]0x201D
    $old = ", substituting this page for the " + $openQuote + "standard form certificate" + $closeQuote + " those terms refer to:"
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Text = "."
    }
}

# ---------------------------------------------------------------------
# 2. Flip "overflow punctuation" off for every paragraph style that
#    currently has it on (it stays on only for styles that never had it
#    enabled to begin with).
# ---------------------------------------------------------------------
$stylesToFix = @(
    "Normal",
    "Arrowhead List",
    "Block Text",
    "Box List",
    "Bullet List",
    "Chapter Heading",
    "Contents Header",
    "Dashed List",
    "Diamond List",
    "Hand List",
    "Heart List",
    "Implies List",
    "Lower Case List",
    "Numbered List",
    "Plain Text",
    "Section Heading",
    "Square List",
    "Star List",
    "Tick List",
    "Triangle List"
)

foreach ($styleName in $stylesToFix) {
    $s = $d.Styles($styleName)
    if ($s -ne $null) {
        $s.ParagraphFormat.HangingPunctuation = $false
    }
}

# ---------------------------------------------------------------------
# 3. Add the new "ListLabel 7" character style (mirrors the existing
#    ListLabel 1-6 styles: quick style, no overrides).
# ---------------------------------------------------------------------
$newStyle = $d.Styles.Add("ListLabel 7", 2)
$newStyle.QuickStyle = $true
